$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The edit re-orders a block of rows (THIAGO ... MARCO/642.75) to sit right
# before the G3C row, drops the CLOTILDE and GABRIELA rows entirely, and
# changes G3C's balance from 238628.98 to 628.98.
#
# Strategy: stage the rows that need to move through a scratch area far
# below the real data (row 1000+) using Copy/PasteSpecial so the original
# inline-string cell typing (e.g. leading-zero account numbers) is
# preserved exactly as-is (a plain .Value = "005064129" assignment would
# get auto-coerced to the number 5064129 by Excel's type inference).
# -----------------------------------------------------------------------

$xlPasteAll = -4104
$xlShiftUp = -4162
$scratchRow = 1000

# Stage 1: stash the rows we need, in their original order, into the scratch area.
#   scratch row 1000       <- row 2  (G3C)
#   scratch row 1001       <- row 4  (THIAGO)
#   scratch rows 1002-1018 <- rows 6-22 (KELMA ... MARCO/642.75, 17 rows)
$ws.Range("A2:C2").Copy()
$ws.Range("A1000:C1000").PasteSpecial($xlPasteAll)

$ws.Range("A4:C4").Copy()
$ws.Range("A1001:C1001").PasteSpecial($xlPasteAll)

$ws.Range("A6:C22").Copy()
$ws.Range("A1002:C1018").PasteSpecial($xlPasteAll)

$excel.CutCopyMode = 0

# Stage 2: update the stashed G3C balance.
$ws.Cells.Item($scratchRow, 3).Value = 628.98

# Stage 3: write the staged rows back into their final destinations.
#   row 2       <- scratch row 1001 (THIAGO)
#   rows 3-19   <- scratch rows 1002-1018 (KELMA ... MARCO/642.75)
#   row 20      <- scratch row 1000 (G3C, now 628.98)
$ws.Range("A1001:C1001").Copy()
$ws.Range("A2:C2").PasteSpecial($xlPasteAll)

$ws.Range("A1002:C1018").Copy()
$ws.Range("A3:C19").PasteSpecial($xlPasteAll)

$ws.Range("A1000:C1000").Copy()
$ws.Range("A20:C20").PasteSpecial($xlPasteAll)

$excel.CutCopyMode = 0

# Stage 4: the old MARIA(659.01)/MARCO(642.75) rows at 21-22 are now stale
# duplicates (their data already lives at rows 18-19) - delete them so
# everything below shifts up into place.
$ws.Rows("21:22").Delete($xlShiftUp)

# Stage 5: clean up the scratch area.
$ws.Rows("1000:1018").ClearContents()
